$wb = $excel.ActiveWorkbook

# "Generate Report for handoff": refresh the Latest Handoff Datetime for the
# 7df83b2c... file row (row 4) on each locale sheet, now that a new handoff
# has gone out.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-15 02:40:38"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-15 02:40:51"
